$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.373.85'
$ws.Range("E2").Value = '  +2.09%  '
$ws.Range("D3").Value = '3.390.59'
$ws.Range("E3").Value = '  +1.80%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '587.60'
$ws.Range("E5").Value = '  +1.06%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '179.85'
$ws.Range("E6").Value = '  +1.55%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("E8").Value = '  +1.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.195'
$ws.Range("E9").Value = '  +6.19%  '
$ws.Range("E10").Value = '  +1.48%  '
$ws.Range("E11").Value = '  +3.01%  '
$ws.Range("E12").Value = '  +3.16%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '678.67'
$ws.Range("E13").Value = '  -0.63%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '8.62'
$ws.Range("E14").Value = '  +2.40%  '
$ws.Range("D15").Value = '3.933.53'
$ws.Range("E15").Value = '  +1.69%  '
$ws.Range("D16").Value = '69.427.33'
$ws.Range("E16").Value = '  +2.16%  '
$ws.Range("E17").Value = '  +1.87%  '
$ws.Range("D18").Value = '3.386.80'
$ws.Range("E18").Value = '  +1.69%  '
$ws.Range("E19").Value = '  +1.22%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.25'
$ws.Range("E20").Value = '  +1.74%  '
$ws.Range("E21").Value = '  +0.70%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.42'
$ws.Range("E22").Value = '  +0.15%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '17.11'
$ws.Range("E23").Value = '  +0.45%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '103.44'
$ws.Range("E24").Value = '  +4.01%  '
$ws.Range("E25").Value = '  +0.19%  '
$ws.Range("E26").Value = '  +1.30%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.70'
$ws.Range("E27").Value = '  +1.64%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '34.15'
$ws.Range("E28").Value = '  +3.06%  '
$ws.Range("E29").Value = '  +1.49%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '556.33'
$ws.Range("E32").Value = '  -1.80%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.61'
$ws.Range("E33").Value = '  +6.87%  '
$ws.Range("E34").Value = '  +0.87%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '58.08'
$ws.Range("E35").Value = '  +1.50%  '
$ws.Range("E36").Value = '  +0.13%  '
$ws.Range("D37").Value = '3.684.54'
$ws.Range("E37").Value = '  -0.54%  '
$ws.Range("E38").Value = '  +5.17%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '34.97'
$ws.Range("E39").Value = '  +1.19%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.27'
$ws.Range("E40").Value = '  +2.66%  '
$ws.Range("E41").Value = '  +1.35%  '
$ws.Range("E42").Value = '  +3.37%  '
$ws.Range("E43").Value = '  +0.69%  '
$ws.Range("E44").Value = '  +3.83%  '
$ws.Range("E46").Value = '  +0.53%  '
$ws.Range("E47").Value = '  +0.80%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.41'
$ws.Range("E48").Value = '  +5.74%  '
$ws.Range("E49").Value = '  -0.01%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '132.48'
$ws.Range("E50").Value = '  +1.78%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.60'
$ws.Range("E51").Value = '  +2.96%  '
